# adapt ppt to k8s v1.10
#
# Slide 4 ("StatefulSet") shapes move up by 237744 EMU (18.72 pt) to make
# room, and the callout warning about "apps/v1beta2" on clusters < 1.9 is
# removed now that the deck targets Kubernetes v1.10+.
#
# Shape.Top/.Left are Single-precision (points) in the PowerPoint object
# model, so naive "Top = Top - 18.72" accumulates rounding error versus the
# EMU values in the OOXML. PointsForEmu searches for a points value whose
# float32-round-trip-then-EMU-floor reproduces the exact target EMU.

function EmuForPoints($pt) {
    $f = [float]$pt
    return [math]::Floor([double]$f * 12700)
}

function PointsForEmu($targetEmu) {
    $base = $targetEmu / 12700.0
    for ($k = 0; $k -lt 5000; $k++) {
        $candidate = $base + ($k * 0.000001)
        if ((EmuForPoints $candidate) -eq $targetEmu) {
            return $candidate
        }
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$targetEmuByName = @{
    "Rectangle 4"                  = 1575325
    "Rectangle 5"                  = 2258278
    "Rectangle 6"                  = 2941231
    "Rectangle 7"                  = 3624184
    "Rectangle 9"                  = 2956633
    "Rectangle 10"                 = 4974610
    "Rectangle 12"                 = 1590726
    "Rectangle 13"                 = 2273679
    "Speech Bubble: Rectangle 17"  = 1706129
    "Rectangle 19"                 = 5585951
    "Rectangle 20"                 = 4363269
    "Rectangle 21"                 = 3769142
    "Speech Bubble: Rectangle 22"  = 2913207
    "Speech Bubble: Rectangle 23"  = 3976871
    "Speech Bubble: Rectangle 25"  = 4970190
    "Speech Bubble: Rectangle 26"  = 5711792
}

foreach ($name in $targetEmuByName.Keys) {
    $shp = $s.Shapes.Item($name)
    $shp.Top = PointsForEmu $targetEmuByName[$name]
}

# Remove the now-obsolete callout about apps/v1beta2 on clusters < 1.9
$s.Shapes.Item("Speech Bubble: Rectangle 15").Delete()
